$d = $word.ActiveDocument

# Update the date line in the first paragraph.
$d.Content.Find.Execute("2025-04-03 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-04 Friday", 2)

# Update the division-problem answers. Cells are addressed directly by
# (row, column) in the table's object model so that duplicate source
# strings (e.g. "42÷5=8, 2" occurs twice) are each replaced with their
# own distinct target value rather than both being overwritten the same way.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "53÷7=7, 4"
$t.Cell(1, 2).Range.Text  = "55÷9=6, 1"
$t.Cell(1, 3).Range.Text  = "96÷2=48, 0"
$t.Cell(1, 4).Range.Text  = "27÷9=3, 0"
$t.Cell(1, 5).Range.Text  = "59÷5=11, 4"

$t.Cell(5, 1).Range.Text  = "64÷8=8, 0"
$t.Cell(5, 2).Range.Text  = "81÷9=9, 0"
$t.Cell(5, 3).Range.Text  = "77÷9=8, 5"
$t.Cell(5, 4).Range.Text  = "41÷3=13, 2"
$t.Cell(5, 5).Range.Text  = "79÷7=11, 2"

$t.Cell(9, 1).Range.Text  = "26÷6=4, 2"
$t.Cell(9, 2).Range.Text  = "16÷7=2, 2"
$t.Cell(9, 3).Range.Text  = "11÷7=1, 4"
$t.Cell(9, 4).Range.Text  = "44÷4=11, 0"
$t.Cell(9, 5).Range.Text  = "74÷5=14, 4"

$t.Cell(13, 1).Range.Text = "51÷3=17, 0"
$t.Cell(13, 2).Range.Text = "13÷3=4, 1"
$t.Cell(13, 3).Range.Text = "20÷8=2, 4"
$t.Cell(13, 4).Range.Text = "17÷8=2, 1"
$t.Cell(13, 5).Range.Text = "10÷9=1, 1"

$t.Cell(17, 1).Range.Text = "34÷7=4, 6"
$t.Cell(17, 2).Range.Text = "46÷3=15, 1"
$t.Cell(17, 3).Range.Text = "46÷4=11, 2"
$t.Cell(17, 4).Range.Text = "38÷2=19, 0"
$t.Cell(17, 5).Range.Text = "97÷7=13, 6"
